# Scheduled runner update: refresh market-board derived figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) for a
# batch of leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 847.0294
$ws.Range("I6").Value = 99.96429000000001
$ws.Range("J6").Value = 4333.3335
$ws.Range("K6").Value = 299.89287
$ws.Range("L6").Value = 13000.0005
$ws.Range("M6").Value = -187.89287
$ws.Range("N6").Value = -13224.0005
$ws.Range("H8").Value = 220
$ws.Range("I8").Value = 220
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 660
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -521
$ws.Range("H11").Value = 158.57143
$ws.Range("I11").Value = 158.57143
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 158.57143
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -18.57142999999999
$ws.Range("H48").Value = 2400
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 2400
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 7200
$ws.Range("N48").Value = -7784
$ws.Range("H56").Value = 2400
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 2400
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 7200
$ws.Range("N56").Value = -8268
$ws.Range("H58").Value = 1663.4375
$ws.Range("I58").Value = 1121.5
$ws.Range("J58").Value = 2566.6667
$ws.Range("K58").Value = 3364.5
$ws.Range("L58").Value = 7700.000100000001
$ws.Range("M58").Value = -3214.5
$ws.Range("N58").Value = -8000.000100000001
$ws.Range("H76").Value = 3181.818
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3181.818
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -5184
$ws.Range("H100").Value = 2566.6667
$ws.Range("I100").Value = 2350
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2350
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1809
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2746
$ws.Range("N113").Value = -9508
$ws.Range("H137").Value = 2467.151
$ws.Range("I137").Value = 2379.4
$ws.Range("J137").Value = 2637.7778
$ws.Range("K137").Value = 7138.200000000001
$ws.Range("L137").Value = 7913.3334
$ws.Range("M137").Value = -4588.200000000001
$ws.Range("N137").Value = -13013.3334
$ws.Range("H138").Value = 2014.1781
$ws.Range("I138").Value = 1368.9524
$ws.Range("J138").Value = 2888.3547
$ws.Range("K138").Value = 4106.857199999999
$ws.Range("L138").Value = 8665.0641
$ws.Range("M138").Value = 1033.142800000001
$ws.Range("N138").Value = -18945.0641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2505.5
$ws.Range("I2").Value = 2505.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2505.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2392.5
$ws.Range("H32").Value = 572306.8
$ws.Range("I32").Value = 615008.75
$ws.Range("J32").Value = 34262.8
$ws.Range("K32").Value = 615008.75
$ws.Range("L32").Value = 34262.8
$ws.Range("M32").Value = -614721.75
$ws.Range("N32").Value = -34836.8
$ws.Range("H61").Value = 1518.2024
$ws.Range("I61").Value = 1285.7794
$ws.Range("J61").Value = 2506
$ws.Range("K61").Value = 1285.7794
$ws.Range("L61").Value = 2506
$ws.Range("M61").Value = -1073.7794
$ws.Range("N61").Value = -2930
$ws.Range("H110").Value = 2345.1667
$ws.Range("I110").Value = 2414.2
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2414.2
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -369.1999999999998
$ws.Range("N110").Value = -6090
$ws.Range("H116").Value = 2505.5
$ws.Range("I116").Value = 2505.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2505.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -211.5
$ws.Range("H136").Value = 1518.2024
$ws.Range("I136").Value = 1285.7794
$ws.Range("J136").Value = 2506
$ws.Range("K136").Value = 3857.3382
$ws.Range("L136").Value = 7518
$ws.Range("M136").Value = -1307.3382
$ws.Range("N136").Value = -12618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2505.5
$ws.Range("I3").Value = 2505.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2505.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2391.5
$ws.Range("H134").Value = 2083.825
$ws.Range("I134").Value = 1301.4517
$ws.Range("J134").Value = 4778.6665
$ws.Range("K134").Value = 3904.3551
$ws.Range("L134").Value = 14335.9995
$ws.Range("M134").Value = -1369.3551
$ws.Range("N134").Value = -19405.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6781.6875
$ws.Range("I31").Value = 1301.7727
$ws.Range("J31").Value = 11418.538
$ws.Range("K31").Value = 1301.7727
$ws.Range("L31").Value = 11418.538
$ws.Range("M31").Value = -1006.7727
$ws.Range("N31").Value = -12008.538
$ws.Range("H34").Value = 6781.6875
$ws.Range("I34").Value = 1301.7727
$ws.Range("J34").Value = 11418.538
$ws.Range("K34").Value = 1301.7727
$ws.Range("L34").Value = 11418.538
$ws.Range("M34").Value = -1099.7727
$ws.Range("N34").Value = -11822.538
$ws.Range("H109").Value = 29800
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 29800
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 29800
$ws.Range("N109").Value = -31880
$ws.Range("H122").Value = 1986.909
$ws.Range("I122").Value = 1556
$ws.Range("J122").Value = 2030
$ws.Range("K122").Value = 4668
$ws.Range("L122").Value = 6090
$ws.Range("M122").Value = -2218
$ws.Range("N122").Value = -10990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 45475.05
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 45475.05
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 136425.15
$ws.Range("N9").Value = -136873.15

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 23470
$ws.Range("I12").Value = 203
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 203
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -63
$ws.Range("N12").Value = -70284
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H52").Value = 14000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 14000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 14000
$ws.Range("N52").Value = -14518
$ws.Range("H80").Value = 84901000
$ws.Range("I80").Value = 127250750
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 127250750
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -127249752
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 84901000
$ws.Range("I83").Value = 127250750
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 636253750
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -636248758
$ws.Range("N83").Value = -1017484
$ws.Range("H102").Value = 1622.1936
$ws.Range("I102").Value = 1366.2084
$ws.Range("J102").Value = 2499.8572
$ws.Range("K102").Value = 1366.2084
$ws.Range("L102").Value = 2499.8572
$ws.Range("M102").Value = 255.7916
$ws.Range("N102").Value = -5743.8572
$ws.Range("H132").Value = 2225
$ws.Range("I132").Value = 2001.9791
$ws.Range("J132").Value = 2788.4211
$ws.Range("K132").Value = 6005.9373
$ws.Range("L132").Value = 8365.263300000001
$ws.Range("M132").Value = -3475.9373
$ws.Range("N132").Value = -13425.2633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 27076.875
$ws.Range("I17").Value = 1320
$ws.Range("J17").Value = 70005
$ws.Range("K17").Value = 1320
$ws.Range("L17").Value = 70005
$ws.Range("M17").Value = -1148
$ws.Range("N17").Value = -70349
$ws.Range("H96").Value = 5125.7144
$ws.Range("I96").Value = 2695
$ws.Range("J96").Value = 8366.666999999999
$ws.Range("K96").Value = 2695
$ws.Range("L96").Value = 8366.666999999999
$ws.Range("M96").Value = -1322
$ws.Range("N96").Value = -11112.667
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
